$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: row, Fecha(D), Calidad(I), Volumen(J), PrecioMin(K), PrecioMax(L), PrecioPromPond(M), Origen(O)
$data = @(
    ,@(371, 44476, "Primera", 11600, 500, 600, 549, "Región Metropolitana")
    ,@(372, 44476, "Primera", 7200, 600, 700, 653, "Región de Coquimbo")
    ,@(373, 44476, "Segunda", 2800, 300, 300, 300, "Región Metropolitana")
    ,@(374, 44386, "Primera", 13700, 500, 600, 550, "Región Metropolitana")
    ,@(375, 44386, "Segunda", 3200, 400, 400, 400, "Región Metropolitana")
    ,@(376, 44384, "Primera", 27300, 450, 600, 508, "Región Metropolitana")
    ,@(377, 44384, "Segunda", 9700, 300, 400, 354, "Región Metropolitana")
    ,@(378, 44263, "Primera", 5500, 800, 800, 800, "Región Metropolitana")
    ,@(379, 44263, "Segunda", 2200, 600, 600, 600, "Región Metropolitana")
    ,@(380, 44363, "Primera", 16700, 400, 500, 441, "Región Metropolitana")
    ,@(381, 44363, "Segunda", 4800, 350, 400, 376, "Región Metropolitana")
    ,@(382, 44306, "Primera", 15600, 400, 500, 430, "Región Metropolitana")
    ,@(383, 44306, "Segunda", 5200, 250, 300, 273, "Región Metropolitana")
    ,@(384, 44369, "Primera", 28300, 350, 500, 423, "Región Metropolitana")
    ,@(385, 44369, "Segunda", 5500, 250, 300, 275, "Región Metropolitana")
    ,@(386, 44172, "Primera", 16000, 700, 800, 741, "Región Metropolitana")
    ,@(387, 44172, "Segunda", 4500, 500, 500, 500, "Región Metropolitana")
    ,@(388, 44301, "Primera", 33300, 400, 600, 518, "Región Metropolitana")
    ,@(389, 44301, "Segunda", 9500, 300, 400, 374, "Región Metropolitana")
    ,@(390, 44357, "Primera", 20400, 300, 450, 369, "Región Metropolitana")
    ,@(391, 44357, "Segunda", 6000, 250, 300, 277, "Región Metropolitana")
    ,@(392, 44328, "Primera", 9100, 400, 500, 436, "Región Metropolitana")
    ,@(393, 44328, "Segunda", 4700, 250, 300, 279, "Región Metropolitana")
    ,@(394, 44328, "Segunda", 300, 500, 500, 500, "Región de O'Higgins")
    ,@(395, 44321, "Primera", 8000, 450, 500, 472, "Región Metropolitana")
    ,@(396, 44321, "Segunda", 6500, 350, 400, 373, "Región Metropolitana")
    ,@(397, 44223, "Primera", 10300, 1000, 1200, 1091, "Región Metropolitana")
    ,@(398, 44223, "Segunda", 11000, 700, 800, 741, "Región Metropolitana")
    ,@(399, 44298, "Primera", 25200, 500, 600, 549, "Región Metropolitana")
    ,@(400, 44298, "Segunda", 8400, 350, 400, 377, "Región Metropolitana")
    ,@(401, 44397, "Primera", 13800, 500, 600, 549, "Región Metropolitana")
    ,@(402, 44397, "Primera", 2300, 700, 700, 700, "Región de O'Higgins")
    ,@(403, 44397, "Segunda", 3600, 500, 500, 500, "Región Metropolitana")
    ,@(404, 44397, "Segunda", 1200, 500, 500, 500, "Región de O'Higgins")
    ,@(405, 44333, "Primera", 9000, 500, 600, 536, "Región Metropolitana")
    ,@(406, 44333, "Primera", 2800, 500, 500, 500, "Región de O'Higgins")
    ,@(407, 44333, "Segunda", 3400, 400, 400, 400, "Región Metropolitana")
    ,@(408, 44314, "Primera", 11300, 350, 400, 379, "Región Metropolitana")
    ,@(409, 44314, "Segunda", 3800, 250, 250, 250, "Región Metropolitana")
    ,@(410, 44392, "Primera", 19400, 550, 700, 609, "Región Metropolitana")
    ,@(411, 44392, "Segunda", 5700, 450, 500, 473, "Región Metropolitana")
    ,@(412, 44425, "Primera", 14600, 500, 600, 542, "Región Metropolitana")
    ,@(413, 44425, "Segunda", 5400, 500, 500, 500, "Región Metropolitana")
    ,@(414, 44390, "Primera", 10200, 500, 600, 566, "Región Metropolitana")
    ,@(415, 44390, "Segunda", 2900, 400, 400, 400, "Región Metropolitana")
    ,@(416, 44187, "Primera", 7000, 600, 700, 650, "Región Metropolitana")
    ,@(417, 44187, "Segunda", 2500, 500, 500, 500, "Región Metropolitana")
    ,@(418, 44466, "Primera", 8400, 450, 500, 473, "Región Metropolitana")
    ,@(419, 44466, "Primera", 6000, 500, 600, 558, "Región de O'Higgins")
    ,@(420, 44466, "Segunda", 3300, 350, 350, 350, "Región Metropolitana")
    ,@(421, 44466, "Segunda", 1500, 400, 400, 400, "Región de O'Higgins")
    ,@(422, 44270, "Primera", 12000, 650, 900, 759, "Región Metropolitana")
    ,@(423, 44270, "Segunda", 4700, 450, 700, 583, "Región Metropolitana")
    ,@(424, 44250, "Primera", 6500, 1000, 1100, 1049, "Región Metropolitana")
    ,@(425, 44250, "Segunda", 2800, 800, 800, 800, "Región Metropolitana")
    ,@(426, 44438, "Primera", 2700, 600, 600, 600, "Provincia de Chacabuco")
    ,@(427, 44438, "Segunda", 2100, 400, 400, 400, "Provincia de Chacabuco")
    ,@(428, 44438, "Tercera", 3300, 550, 550, 550, "Provincia de Chacabuco")
    ,@(429, 44201, "Primera", 5200, 700, 800, 750, "Región Metropolitana")
    ,@(430, 44201, "Segunda", 2500, 600, 600, 600, "Región Metropolitana")
    ,@(431, 44461, "Primera", 12500, 500, 600, 554, "Región Metropolitana")
    ,@(432, 44461, "Segunda", 2600, 500, 500, 500, "Región Metropolitana")
    ,@(433, 44193, "Primera", 6100, 700, 800, 750, "Región Metropolitana")
    ,@(434, 44193, "Segunda", 2500, 600, 600, 600, "Región Metropolitana")
    ,@(435, 44286, "Primera", 7900, 700, 800, 757, "Región Metropolitana")
    ,@(436, 44286, "Primera", 6800, 700, 800, 747, "Región de O'Higgins")
    ,@(437, 44286, "Segunda", 2500, 500, 500, 500, "Región Metropolitana")
    ,@(438, 44286, "Segunda", 2300, 600, 600, 600, "Región de O'Higgins")
    ,@(439, 44389, "Primera", 3500, 700, 700, 700, "Región Metropolitana")
    ,@(440, 44389, "Segunda", 2800, 600, 600, 600, "Región Metropolitana")
    ,@(441, 44312, "Primera", 14800, 350, 500, 414, "Región Metropolitana")
    ,@(442, 44312, "Segunda", 5200, 250, 300, 277, "Región Metropolitana")
    ,@(443, 44326, "Primera", 8700, 400, 600, 511, "Región Metropolitana")
    ,@(444, 44326, "Segunda", 4000, 250, 400, 355, "Región Metropolitana")
    ,@(445, 44432, "Primera", 16300, 400, 500, 442, "Región Metropolitana")
    ,@(446, 44432, "Segunda", 6100, 350, 400, 378, "Región Metropolitana")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]   # D Fecha
    $ws.Cells.Item($r, 9).Value = $row[2]   # I Calidad
    $ws.Cells.Item($r, 10).Value = $row[3]  # J Volumen
    $ws.Cells.Item($r, 11).Value = $row[4]  # K Precio minimo
    $ws.Cells.Item($r, 12).Value = $row[5]  # L Precio maximo
    $ws.Cells.Item($r, 13).Value = $row[6]  # M Precio promedio ponderado
    $ws.Cells.Item($r, 15).Value = $row[7]  # O Origen
    $ws.Cells.Item($r, 16).Value = $row[6]  # P Precio $/Kg (mirrors M)
}

# New rows 444-446 need all constant columns filled in too
# (rows 371..443 already existed with these constant columns populated)
$dateNumberFormat = $ws.Cells.Item(370, 4).NumberFormat
foreach ($r in 444..446) {
    $ws.Cells.Item($r, 1).Value = 6                                   # A Mercado ID
    $ws.Cells.Item($r, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"  # B Mercado
    $ws.Cells.Item($r, 3).Value = "Metropolitana"                     # C Region
    $ws.Cells.Item($r, 4).NumberFormat = $dateNumberFormat            # D style (date fmt), value set above
    $ws.Cells.Item($r, 5).Value = 13                                  # E Codreg
    $ws.Cells.Item($r, 6).Value = 100112023                           # F Categoria ID
    $ws.Cells.Item($r, 7).Value = "Brócoli"                           # G Categoria
    $ws.Cells.Item($r, 8).Value = "Sin especificar"                   # H Variedad
    $ws.Cells.Item($r, 14).Value = "`$/unidad"                        # N Unidad de comercializacion
    $ws.Cells.Item($r, 17).Value = 1                                  # Q Kg o Unidades
    $ws.Cells.Item($r, 18).Value = "Hortaliza"                        # R Clasificacion
}

Write-Output "done"
